$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.65714501444922
$ws.Range("C2").Value = 5.35832527971197
$ws.Range("D2").Value = 6.030266841174345
$ws.Range("E2").Value = 10.6149482203283
$ws.Range("G2").Value = 52.35904211565074
$ws.Range("H2").Value = 19.94296560989447
$ws.Range("I2").Value = 31.56467645502768
$ws.Range("K2").Value = 11.6290452005049
$ws.Range("L2").Value = 10.17709072562915

$ws.Range("B3").Value = 13.49898012448792
$ws.Range("C3").Value = 5.202703433163136
$ws.Range("D3").Value = 5.919629805366032
$ws.Range("E3").Value = 10.62054612710525
$ws.Range("G3").Value = 52.05360397563903
$ws.Range("H3").Value = 19.93632045796619
$ws.Range("I3").Value = 31.53071451958706
$ws.Range("K3").Value = 11.5218807811543
$ws.Range("L3").Value = 10.1685448459058

$ws.Range("B4").Value = 13.40508190964185
$ws.Range("C4").Value = 5.103452716969935
$ws.Range("D4").Value = 5.852492666321898
$ws.Range("E4").Value = 10.62559304327286
$ws.Range("G4").Value = 51.87631122658168
$ws.Range("H4").Value = 19.93522807215466
$ws.Range("I4").Value = 31.51481978652997
$ws.Range("K4").Value = 11.45896476709947
$ws.Range("L4").Value = 10.16521584036634

$ws.Range("B5").Value = 13.36767466087269
$ws.Range("C5").Value = 5.062104900677698
$ws.Range("D5").Value = 5.825373893573653
$ws.Range("E5").Value = 10.62805481982173
$ws.Range("G5").Value = 51.80668777963565
$ws.Range("H5").Value = 19.93553455152819
$ws.Range("I5").Value = 31.50959116799577
$ws.Range("K5").Value = 11.43408006371554
$ws.Range("L5").Value = 10.16434296052382

$ws.Range("B6").Value = 13.36151635835295
$ws.Range("C6").Value = 5.055185485775334
$ws.Range("D6").Value = 5.82088652645276
$ws.Range("E6").Value = 10.62848807142011
$ws.Range("G6").Value = 51.79528670727586
$ws.Range("H6").Value = 19.9356308428366
$ws.Range("I6").Value = 31.50879841369827
$ws.Range("K6").Value = 11.4299943407576
$ws.Range("L6").Value = 10.16422726927484

$ws.Range("B7").Value = 13.40457388966091
$ws.Range("C7").Value = 5.102898701005539
$ws.Range("D7").Value = 5.852125907321283
$ws.Range("E7").Value = 10.62562460299352
$ws.Range("G7").Value = 51.87536157271943
$ws.Range("H7").Value = 19.93522916194615
$ws.Range("I7").Value = 31.51474421410122
$ws.Range("K7").Value = 11.45862607282628
$ws.Range("L7").Value = 10.16520210825398

$ws.Range("B8").Value = 13.60197193263899
$ws.Range("C8").Value = 5.305450211580684
$ws.Range("D8").Value = 5.991978511090287
$ws.Range("E8").Value = 10.61654449345137
$ws.Range("G8").Value = 52.25163022868752
$ws.Range("H8").Value = 19.94005453871538
$ws.Range("I8").Value = 31.55193866762967
$ws.Range("K8").Value = 11.59151434206853
$ws.Range("L8").Value = 10.17374688754072

$ws.Range("B9").Value = 14.01237748201029
$ws.Range("C9").Value = 5.672227225454589
$ws.Range("D9").Value = 6.27076892633707
$ws.Range("E9").Value = 10.61149387591887
$ws.Range("G9").Value = 53.06850320775963
$ws.Range("H9").Value = 19.97319902222403
$ws.Range("I9").Value = 31.6641087426533
$ws.Range("K9").Value = 11.87364061981226
$ws.Range("L9").Value = 10.20565205559023

$ws.Range("B10").Value = 14.32507835937608
$ws.Range("C10").Value = 5.921876908221743
$ws.Range("D10").Value = 6.476046081043905
$ws.Range("E10").Value = 10.61553116425293
$ws.Range("G10").Value = 53.71349427989919
$ws.Range("H10").Value = 20.01193990180071
$ws.Range("I10").Value = 31.77025506416917
$ws.Range("K10").Value = 12.0921615807107
$ws.Range("L10").Value = 10.23821650297603

$ws.Range("B11").Value = 14.46910564432527
$ws.Range("C11").Value = 6.030921715208457
$ws.Range("D11").Value = 6.569059818393475
$ws.Range("E11").Value = 10.61904194824692
$ws.Range("G11").Value = 54.0158712108013
$ws.Range("H11").Value = 20.03267200983704
$ws.Range("I11").Value = 31.8236437996157
$ws.Range("K11").Value = 12.19359586005945
$ws.Range("L11").Value = 10.25498235811142

$ws.Range("B12").Value = 14.52384609028485
$ws.Range("C12").Value = 6.071547883132284
$ws.Range("D12").Value = 6.604191646011623
$ws.Range("E12").Value = 10.62061120036375
$ws.Range("G12").Value = 54.13159250035513
$ws.Range("H12").Value = 20.04096766415697
$ws.Range("I12").Value = 31.84458845463424
$ws.Range("K12").Value = 12.23226235243428
$ws.Range("L12").Value = 10.26160884536168

$ws.Range("B13").Value = 14.51204871364459
$ws.Range("C13").Value = 6.062828205925884
$ws.Range("D13").Value = 6.596629980431652
$ws.Range("E13").Value = 10.62026258544042
$ws.Range("G13").Value = 54.10661686075721
$ws.Range("H13").Value = 20.03916129972832
$ws.Range("I13").Value = 31.84004539929736
$ws.Range("K13").Value = 12.2239240422651
$ws.Range("L13").Value = 10.26016941774001

$ws.Range("B14").Value = 14.47360547072954
$ws.Range("C14").Value = 6.034277504309721
$ws.Range("D14").Value = 6.57195215666557
$ws.Range("E14").Value = 10.61916625163218
$ws.Range("G14").Value = 54.0253676860165
$ws.Range("H14").Value = 20.03334559641221
$ws.Range("I14").Value = 31.82535238697312
$ws.Range("K14").Value = 12.19677206150201
$ws.Range("L14").Value = 10.25552197797817

$ws.Range("B15").Value = 14.45008234900993
$ws.Range("C15").Value = 6.016702091003253
$ws.Range("D15").Value = 6.556823412618231
$ws.Range("E15").Value = 10.61852591339271
$ws.Range("G15").Value = 53.97575663003403
$ws.Range("H15").Value = 20.02984117315766
$ws.Range("I15").Value = 31.81644703788452
$ws.Range("K15").Value = 12.18017289406371
$ws.Range("L15").Value = 10.25271134483305

$ws.Range("B16").Value = 14.31569688066552
$ws.Range("C16").Value = 5.914658126009953
$ws.Range("D16").Value = 6.469956761688541
$ws.Range("E16").Value = 10.61533533680403
$ws.Range("G16").Value = 53.69390760251353
$ws.Range("H16").Value = 20.0106473868364
$ws.Range("I16").Value = 31.76686814235833
$ws.Range("K16").Value = 12.08557033863614
$ws.Range("L16").Value = 10.23715983897753

$ws.Range("B17").Value = 14.2336716715648
$ws.Range("C17").Value = 5.850886908764855
$ws.Range("D17").Value = 6.416545394131918
$ws.Range("E17").Value = 10.61380618564781
$ws.Range("G17").Value = 53.52324882936595
$ws.Range("H17").Value = 19.99966739518092
$ws.Range("I17").Value = 31.7377557327321
$ws.Range("K17").Value = 12.028028530075
$ws.Range("L17").Value = 10.22811737209134

$ws.Range("B18").Value = 14.18666358282044
$ws.Range("C18").Value = 5.813782606925113
$ws.Range("D18").Value = 6.385792104398795
$ws.Range("E18").Value = 10.61308430709408
$ws.Range("G18").Value = 53.42593892534357
$ws.Range("H18").Value = 19.99364475036627
$ws.Range("I18").Value = 31.72149171050128
$ws.Range("K18").Value = 11.99512506575438
$ws.Range("L18").Value = 10.22310033101211

$ws.Range("B19").Value = 14.17077845259423
$ws.Range("C19").Value = 5.80114730429708
$ws.Range("D19").Value = 6.375375162607923
$ws.Range("E19").Value = 10.61286699474998
$ws.Range("G19").Value = 53.39313934145085
$ws.Range("H19").Value = 19.99165592799553
$ws.Range("I19").Value = 31.71606772339879
$ws.Range("K19").Value = 11.98401877377174
$ws.Range("L19").Value = 10.22143332959545

$ws.Range("B20").Value = 14.24238613645504
$ws.Range("C20").Value = 5.857719552713029
$ws.Range("D20").Value = 6.422234764990659
$ws.Range("E20").Value = 10.61395265686232
$ws.Range("G20").Value = 53.54132844153118
$ws.Range("H20").Value = 20.00080594946781
$ws.Range("I20").Value = 31.74080509376063
$ws.Range("K20").Value = 12.03413423362634
$ws.Range("L20").Value = 10.22906094181382

$ws.Range("B21").Value = 14.48489218016162
$ws.Range("C21").Value = 6.04268175410864
$ws.Range("D21").Value = 6.579203377147071
$ws.Range("E21").Value = 10.61948177171653
$ws.Range("G21").Value = 54.04920007395022
$ws.Range("H21").Value = 20.03504175521354
$ws.Range("I21").Value = 31.82964838981852
$ws.Range("K21").Value = 12.20474060948947
$ws.Range("L21").Value = 10.25687953517565

$ws.Range("B22").Value = 14.6445271302591
$ws.Range("C22").Value = 6.159672698459537
$ws.Range("D22").Value = 6.681248057929018
$ws.Range("E22").Value = 10.62449254253711
$ws.Range("G22").Value = 54.38818635502457
$ws.Range("H22").Value = 20.06000853956349
$ws.Range("I22").Value = 31.89194964399967
$ws.Range("K22").Value = 12.31771441193938
$ws.Range("L22").Value = 10.27667715127544

$ws.Range("B23").Value = 14.55924040958974
$ws.Range("C23").Value = 6.097593547172473
$ws.Range("D23").Value = 6.626846551756264
$ws.Range("E23").Value = 10.62169070564847
$ws.Range("G23").Value = 54.20664107081369
$ws.Range("H23").Value = 20.04644696142983
$ws.Range("I23").Value = 31.85831284418722
$ws.Range("K23").Value = 12.25729535644679
$ws.Range("L23").Value = 10.26596397094006

$ws.Range("B24").Value = 14.23844585759615
$ws.Range("C24").Value = 5.854631887995963
$ws.Range("D24").Value = 6.419662742057514
$ws.Range("E24").Value = 10.61388594734512
$ws.Range("G24").Value = 53.53315213856958
$ws.Range("H24").Value = 20.0002903058022
$ws.Range("I24").Value = 31.73942500341516
$ws.Range("K24").Value = 12.03137328755675
$ws.Range("L24").Value = 10.22863378806751

$ws.Range("B25").Value = 13.89916707205997
$ws.Range("C25").Value = 5.576402105090348
$ws.Range("D25").Value = 6.195105707013814
$ws.Range("E25").Value = 10.61149707226665
$ws.Range("G25").Value = 52.83938840873979
$ws.Range("H25").Value = 19.96170082492816
$ws.Range("I25").Value = 31.62957949552369
$ws.Range("K25").Value = 11.79520361055847
$ws.Range("L25").Value = 10.19540940689224
